# Weekly fruit/vegetable price update: the 4 "Caqui" sample rows (2, 3, 4, 7)
# get their Fecha/Volumen/Precio/Unidad/Origen/Precio-Kg/Kg-unidad values
# rotated to new observations, per the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44344
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 13500
$ws.Range("Q2").Value = "$/caja 18 kilos granel"
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 750
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44330
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 15500
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value = 861

# Row 4
$ws.Range("D4").Value = 44698
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 16500
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 917

# Row 7
$ws.Range("D7").Value = 44334
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 12500
$ws.Range("Q7").Value = "$/caja 12 kilos empedrada"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1042
$ws.Range("T7").Value = 12
